$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.063857049188834
$ws.Cells.Item(2, 4).Value2 = 1.062743240390809
$ws.Cells.Item(2, 5).Value2 = 1.0686520845984
$ws.Cells.Item(2, 6).Value2 = 1.078065579570893
$ws.Cells.Item(2, 9).Value2 = 1.05206115945295
$ws.Cells.Item(2, 10).Value2 = 1.068820272934943
$ws.Cells.Item(2, 11).Value2 = 1.065463760880149
$ws.Cells.Item(2, 12).Value2 = 1.071356694158912
$ws.Cells.Item(2, 13).Value2 = 1.080745230416667
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.065093264452941
$ws.Cells.Item(3, 4).Value2 = 1.063711229497851
$ws.Cells.Item(3, 5).Value2 = 1.069763634194439
$ws.Cells.Item(3, 6).Value2 = 1.079280221981993
$ws.Cells.Item(3, 9).Value2 = 1.052454004787267
$ws.Cells.Item(3, 10).Value2 = 1.069710114300826
$ws.Cells.Item(3, 11).Value2 = 1.066246240683857
$ws.Cells.Item(3, 12).Value2 = 1.072283528993061
$ws.Cells.Item(3, 13).Value2 = 1.081776717094959
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.065892923662332
$ws.Cells.Item(4, 4).Value2 = 1.064337270394392
$ws.Cells.Item(4, 5).Value2 = 1.070482917961386
$ws.Cells.Item(4, 6).Value2 = 1.08006639377219
$ws.Cells.Item(4, 9).Value2 = 1.052706863294401
$ws.Cells.Item(4, 10).Value2 = 1.070285102132902
$ws.Cells.Item(4, 11).Value2 = 1.066751623846446
$ws.Cells.Item(4, 12).Value2 = 1.072882696669001
$ws.Cells.Item(4, 13).Value2 = 1.082443799659464
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.066229041488041
$ws.Cells.Item(5, 4).Value2 = 1.064600384040476
$ws.Cells.Item(5, 5).Value2 = 1.070785314909472
$ws.Cells.Item(5, 6).Value2 = 1.080396953763106
$ws.Cells.Item(5, 9).Value2 = 1.052812845102037
$ws.Cells.Item(5, 10).Value2 = 1.070526636950645
$ws.Cells.Item(5, 11).Value2 = 1.066963864412702
$ws.Cells.Item(5, 12).Value2 = 1.073134454712385
$ws.Cells.Item(5, 13).Value2 = 1.082724156683313
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.066285473756093
$ws.Cells.Item(6, 4).Value2 = 1.064644557670262
$ws.Cells.Item(6, 5).Value2 = 1.070836089312904
$ws.Cells.Item(6, 6).Value2 = 1.080452459411556
$ws.Cells.Item(6, 9).Value2 = 1.052830621171062
$ws.Cells.Item(6, 10).Value2 = 1.070567180599296
$ws.Cells.Item(6, 11).Value2 = 1.06699948751801
$ws.Cells.Item(6, 12).Value2 = 1.073176718286175
$ws.Cells.Item(6, 13).Value2 = 1.082771224944406
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.065897415114027
$ws.Cells.Item(7, 4).Value2 = 1.064340786421202
$ws.Cells.Item(7, 5).Value2 = 1.070486958563168
$ws.Cells.Item(7, 6).Value2 = 1.080070810519458
$ws.Cells.Item(7, 9).Value2 = 1.052708280684802
$ws.Cells.Item(7, 10).Value2 = 1.070288330278262
$ws.Cells.Item(7, 11).Value2 = 1.066754460688379
$ws.Cells.Item(7, 12).Value2 = 1.072886061189701
$ws.Cells.Item(7, 13).Value2 = 1.082447546134099
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.064274887918893
$ws.Cells.Item(8, 4).Value2 = 1.063070442168659
$ws.Cells.Item(8, 5).Value2 = 1.069027731099451
$ws.Cells.Item(8, 6).Value2 = 1.078476029609322
$ws.Cells.Item(8, 9).Value2 = 1.052194200905632
$ws.Cells.Item(8, 10).Value2 = 1.069121165172783
$ws.Cells.Item(8, 11).Value2 = 1.065728397495771
$ws.Cells.Item(8, 12).Value2 = 1.07167003824608
$ws.Cells.Item(8, 13).Value2 = 1.081093901131599
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.061413743593206
$ws.Cells.Item(9, 4).Value2 = 1.060829495441938
$ws.Cells.Item(9, 5).Value2 = 1.066456612621051
$ws.Cells.Item(9, 6).Value2 = 1.075667418608601
$ws.Cells.Item(9, 9).Value2 = 1.051278047836727
$ws.Cells.Item(9, 10).Value2 = 1.067058296391802
$ws.Cells.Item(9, 11).Value2 = 1.063913148223366
$ws.Cells.Item(9, 12).Value2 = 1.069522936318202
$ws.Cells.Item(9, 13).Value2 = 1.078705803555414
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.059504807249836
$ws.Cells.Item(10, 4).Value2 = 1.059333830171813
$ws.Cells.Item(10, 5).Value2 = 1.064742605969384
$ws.Cells.Item(10, 6).Value2 = 1.073795991791327
$ws.Cells.Item(10, 9).Value2 = 1.050660330551296
$ws.Cells.Item(10, 10).Value2 = 1.065678819800986
$ws.Cells.Item(10, 11).Value2 = 1.062698080156298
$ws.Cells.Item(10, 12).Value2 = 1.068088558832393
$ws.Cells.Item(10, 13).Value2 = 1.077111770974686
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.058677829183443
$ws.Cells.Item(11, 4).Value2 = 1.058685772652103
$ws.Cells.Item(11, 5).Value2 = 1.064000420694308
$ws.Cells.Item(11, 6).Value2 = 1.07298585437603
$ws.Cells.Item(11, 9).Value2 = 1.050391195233523
$ws.Cells.Item(11, 10).Value2 = 1.065080469426528
$ws.Cells.Item(11, 11).Value2 = 1.06217076458821
$ws.Cells.Item(11, 12).Value2 = 1.067466734095875
$ws.Cells.Item(11, 13).Value2 = 1.076421051263113
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.058370590515273
$ws.Cells.Item(12, 4).Value2 = 1.058444990092801
$ws.Cells.Item(12, 5).Value2 = 1.063724736970102
$ws.Cells.Item(12, 6).Value2 = 1.072684961777557
$ws.Cells.Item(12, 9).Value2 = 1.050290976223825
$ws.Cells.Item(12, 10).Value2 = 1.064858059353596
$ws.Cells.Item(12, 11).Value2 = 1.061974716978239
$ws.Cells.Item(12, 12).Value2 = 1.067235649451026
$ws.Cells.Item(12, 13).Value2 = 1.076164411365345
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.058436497110218
$ws.Cells.Item(13, 4).Value2 = 1.058496641738697
$ws.Cells.Item(13, 5).Value2 = 1.063783872230966
$ws.Cells.Item(13, 6).Value2 = 1.072749503002633
$ws.Cells.Item(13, 9).Value2 = 1.050312484875556
$ws.Cells.Item(13, 10).Value2 = 1.06490577416006
$ws.Cells.Item(13, 11).Value2 = 1.062016777984339
$ws.Cells.Item(13, 12).Value2 = 1.067285222951919
$ws.Cells.Item(13, 13).Value2 = 1.07621946495498
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.058652434023605
$ws.Cells.Item(14, 4).Value2 = 1.058665870827348
$ws.Cells.Item(14, 5).Value2 = 1.0639776326641
$ws.Cells.Item(14, 6).Value2 = 1.072960981929705
$ws.Cells.Item(14, 9).Value2 = 1.050382916207072
$ws.Cells.Item(14, 10).Value2 = 1.065062088134636
$ws.Cells.Item(14, 11).Value2 = 1.06215456289382
$ws.Cells.Item(14, 12).Value2 = 1.06744763484582
$ws.Cells.Item(14, 13).Value2 = 1.076399838886476
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.058785471560916
$ws.Cells.Item(15, 4).Value2 = 1.058770129781686
$ws.Cells.Item(15, 5).Value2 = 1.064097014404407
$ws.Cells.Item(15, 6).Value2 = 1.073091284795619
$ws.Cells.Item(15, 9).Value2 = 1.050426278100382
$ws.Cells.Item(15, 10).Value2 = 1.065158377598268
$ws.Cells.Item(15, 11).Value2 = 1.062239432946161
$ws.Cells.Item(15, 12).Value2 = 1.067547687384048
$ws.Cells.Item(15, 13).Value2 = 1.076510963116101
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.059559682486749
$ws.Cells.Item(16, 4).Value2 = 1.059376830554194
$ws.Cells.Item(16, 5).Value2 = 1.064791862025169
$ws.Cells.Item(16, 6).Value2 = 1.073849762048636
$ws.Cells.Item(16, 9).Value2 = 1.050678157138525
$ws.Cells.Item(16, 10).Value2 = 1.065718508545365
$ws.Cells.Item(16, 11).Value2 = 1.062733051309648
$ws.Cells.Item(16, 12).Value2 = 1.068129811778961
$ws.Cells.Item(16, 13).Value2 = 1.077157601228065
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.060045216616109
$ws.Cells.Item(17, 4).Value2 = 1.059757283292209
$ws.Cells.Item(17, 5).Value2 = 1.065227717947015
$ws.Cells.Item(17, 6).Value2 = 1.074325587694898
$ws.Cells.Item(17, 9).Value2 = 1.050835709285125
$ws.Cells.Item(17, 10).Value2 = 1.066069587511152
$ws.Cells.Item(17, 11).Value2 = 1.063042367347255
$ws.Cells.Item(17, 12).Value2 = 1.068494766507409
$ws.Cells.Item(17, 13).Value2 = 1.077563086975599
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.060328382431119
$ws.Cells.Item(18, 4).Value2 = 1.059979154063713
$ws.Cells.Item(18, 5).Value2 = 1.065481944572084
$ws.Cells.Item(18, 6).Value2 = 1.0746031484394
$ws.Cells.Item(18, 9).Value2 = 1.050927446710176
$ws.Cells.Item(18, 10).Value2 = 1.066274266868013
$ws.Cells.Item(18, 11).Value2 = 1.063222672103732
$ws.Cells.Item(18, 12).Value2 = 1.068707568118529
$ws.Cells.Item(18, 13).Value2 = 1.077799552496662
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.060424928287102
$ws.Cells.Item(19, 4).Value2 = 1.060054799370613
$ws.Cells.Item(19, 5).Value2 = 1.065568629203464
$ws.Cells.Item(19, 6).Value2 = 1.074697792879599
$ws.Cells.Item(19, 9).Value2 = 1.050958699674368
$ws.Cells.Item(19, 10).Value2 = 1.066344040486406
$ws.Cells.Item(19, 11).Value2 = 1.06328413206196
$ws.Cells.Item(19, 12).Value2 = 1.068780116091314
$ws.Cells.Item(19, 13).Value2 = 1.077880173132099
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.059993127353212
$ws.Cells.Item(20, 4).Value2 = 1.059716468552456
$ws.Cells.Item(20, 5).Value2 = 1.065180954823458
$ws.Cells.Item(20, 6).Value2 = 1.074274534084269
$ws.Cells.Item(20, 9).Value2 = 1.050818821994782
$ws.Cells.Item(20, 10).Value2 = 1.066031930320374
$ws.Cells.Item(20, 11).Value2 = 1.063009192471851
$ws.Cells.Item(20, 12).Value2 = 1.068455617603699
$ws.Cells.Item(20, 13).Value2 = 1.077519587114633
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.058588847723191
$ws.Cells.Item(21, 4).Value2 = 1.058616038887781
$ws.Cells.Item(21, 5).Value2 = 1.063920575144263
$ws.Cells.Item(21, 6).Value2 = 1.072898705885367
$ws.Cells.Item(21, 9).Value2 = 1.050362182847925
$ws.Cells.Item(21, 10).Value2 = 1.065016061887774
$ws.Cells.Item(21, 11).Value2 = 1.062113993628166
$ws.Cells.Item(21, 12).Value2 = 1.067399811675049
$ws.Cells.Item(21, 13).Value2 = 1.076346725367153
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.057705560216798
$ws.Cells.Item(22, 4).Value2 = 1.057923777518498
$ws.Cells.Item(22, 5).Value2 = 1.063128105822981
$ws.Cells.Item(22, 6).Value2 = 1.072033832064472
$ws.Cells.Item(22, 9).Value2 = 1.050073627880074
$ws.Cells.Item(22, 10).Value2 = 1.064376441094569
$ws.Cells.Item(22, 11).Value2 = 1.061550109509574
$ws.Cells.Item(22, 12).Value2 = 1.066735340122136
$ws.Cells.Item(22, 13).Value2 = 1.075608860958362
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.05817384277179
$ws.Cells.Item(23, 4).Value2 = 1.058290794610076
$ws.Cells.Item(23, 5).Value2 = 1.063548211136578
$ws.Cells.Item(23, 6).Value2 = 1.072492303151141
$ws.Cells.Item(23, 9).Value2 = 1.050226733839835
$ws.Cells.Item(23, 10).Value2 = 1.06471560242604
$ws.Cells.Item(23, 11).Value2 = 1.061849133916374
$ws.Cells.Item(23, 12).Value2 = 1.067087650688796
$ws.Cells.Item(23, 13).Value2 = 1.076000059108402
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.060016664366464
$ws.Cells.Item(24, 4).Value2 = 1.05973491110289
$ws.Cells.Item(24, 5).Value2 = 1.065202085068531
$ws.Cells.Item(24, 6).Value2 = 1.074297602952558
$ws.Cells.Item(24, 9).Value2 = 1.050826453129025
$ws.Cells.Item(24, 10).Value2 = 1.06604894629009
$ws.Cells.Item(24, 11).Value2 = 1.063024183122517
$ws.Cells.Item(24, 12).Value2 = 1.068473307525438
$ws.Cells.Item(24, 13).Value2 = 1.077539242976066
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.062153673845778
$ws.Cells.Item(25, 4).Value2 = 1.061409129120433
$ws.Cells.Item(25, 5).Value2 = 1.067121289740179
$ws.Cells.Item(25, 6).Value2 = 1.07639333221923
$ws.Cells.Item(25, 9).Value2 = 1.051516117333635
$ws.Cells.Item(25, 10).Value2 = 1.067592337298529
$ws.Cells.Item(25, 11).Value2 = 1.064383292923893
$ws.Cells.Item(25, 12).Value2 = 1.070078532653068
$ws.Cells.Item(25, 13).Value2 = 1.079323524676657
